$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: fill in the first "missing item" data row ---

# م (sequence number) -> 1
$ws.Range("A7").Value = 1

# الاسم (item name) -> merged C7:G7; also apply Text format (matches N7:O7 below,
# which shares the same underlying style) since both become string-typed data.
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "DRAMENEX 50MG 20 TABS."

# الرصيد الحالي (current balance) -> merged H7:K7, text "1:0"
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "1:0"

# حد الطلب (order limit) -> merged L7:M7, textual "1" but keep the cell's
# existing numeric display format (#,##0.##;"["#,##0.##"]";0) unchanged.
$fmtL7 = $ws.Range("L7").NumberFormat()
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = $fmtL7

# السعر (price) -> merged N7:O7, text "28.00" (shares style with C7:G7)
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "28.00"

# سعر البيع (sell price) -> P7, textual "14.0000" but keep existing numeric
# display format (0.00) unchanged.
$fmtP7 = $ws.Range("P7").NumberFormat()
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "14.0000"
$ws.Range("P7").NumberFormat = $fmtP7

# عدد التعااملات (transaction count) -> Q7, text "0:1"
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "0:1"

# --- Row 8: footer/total cell ---
$ws.Range("N8").Value = 14
